$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.957.42'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '2.238.90'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '242.38'
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").Value = '0.621'
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").Value = '74.18'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -3.79%  '
$ws.Range("D10").Value = '42.11'
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("D11").Value = '0.0950'
$ws.Range("E11").Value = '  -1.57%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").Value = '2.571.82'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("E16").Value = '  -1.66%  '
$ws.Range("D17").Value = '2.240.18'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '41.914.58'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("E19").Value = '  -6.24%  '
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D21").Value = '72.41'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '11.12'
$ws.Range("E22").Value = '  +8.25%  '
$ws.Range("D23").Value = '229.53'
$ws.Range("E23").Value = '  -0.60%  '
$ws.Range("E24").Value = '  -5.69%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").Value = '11.30'
$ws.Range("E26").Value = '  -3.41%  '
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("E28").Value = '  -1.00%  '
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").Value = '167.38'
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").Value = '20.55'
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").Value = '5.58'
$ws.Range("E32").Value = '  -4.45%  '
$ws.Range("D33").Value = '0.0798'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("D34").Value = '30.42'
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("E35").Value = '  -0.53%  '
$ws.Range("E36").Value = '  -6.41%  '
$ws.Range("E37").Value = '  -3.37%  '
$ws.Range("D38").Value = '0.0303'
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("D39").Value = '13.06'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("E40").Value = '  -1.86%  '
$ws.Range("D41").Value = '5.65'
$ws.Range("D42").Value = '64.26'
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("E44").Value = '  -1.41%  '
$ws.Range("D45").Value = '102.99'
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("E48").Value = '  -0.71%  '
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("D51").Value = '2.446.70'
$ws.Range("E51").Value = '  +0.13%  '
